# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the Leve tables on each job sheet, per scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2326.182
$ws.Range("I2").Value = 2358.9
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 2358.9
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -2245.9
$ws.Range("N2").Value = -2225
$ws.Range("H43").Value = 1352.4736
$ws.Range("I43").Value = 1441.7273
$ws.Range("J43").Value = 1229.75
$ws.Range("K43").Value = 1441.7273
$ws.Range("L43").Value = 1229.75
$ws.Range("M43").Value = -1372.7273
$ws.Range("N43").Value = -1367.75
$ws.Range("H58").Value = 423.57144
$ws.Range("I58").Value = 160.83333
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 482.49999
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -332.49999
$ws.Range("N58").Value = -6300
$ws.Range("H88").Value = 2912.1785
$ws.Range("I88").Value = 1205
$ws.Range("J88").Value = 3196.7083
$ws.Range("K88").Value = 1205
$ws.Range("L88").Value = 3196.7083
$ws.Range("M88").Value = -799
$ws.Range("N88").Value = -4008.7083
$ws.Range("H91").Value = 2912.1785
$ws.Range("I91").Value = 1205
$ws.Range("J91").Value = 3196.7083
$ws.Range("K91").Value = 1205
$ws.Range("L91").Value = 3196.7083
$ws.Range("M91").Value = 199
$ws.Range("N91").Value = -6004.7083

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2100
$ws.Range("I35").Value = 1800
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = -1394
$ws.Range("N35").Value = -3812
$ws.Range("H63").Value = 2778.2
$ws.Range("I63").Value = 2778.2
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2778.2
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2092.2
$ws.Range("H66").Value = 2778.2
$ws.Range("I66").Value = 2778.2
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 13891
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -10459
$ws.Range("H88").Value = 1299
$ws.Range("I88").Value = 475.5
$ws.Range("J88").Value = 1665
$ws.Range("K88").Value = 475.5
$ws.Range("L88").Value = 1665
$ws.Range("M88").Value = -69.5
$ws.Range("N88").Value = -2477
$ws.Range("H91").Value = 1299
$ws.Range("I91").Value = 475.5
$ws.Range("J91").Value = 1665
$ws.Range("K91").Value = 475.5
$ws.Range("L91").Value = 1665
$ws.Range("M91").Value = 928.5
$ws.Range("N91").Value = -4473
$ws.Range("H110").Value = 1212864.5
$ws.Range("I110").Value = 1393929.2
$ws.Range("J110").Value = 5766
$ws.Range("K110").Value = 1393929.2
$ws.Range("L110").Value = 5766
$ws.Range("M110").Value = -1391884.2
$ws.Range("N110").Value = -9856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 295000
$ws.Range("I29").Value = 295000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 295000
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -294711
$ws.Range("H36").Value = 5958.75
$ws.Range("I36").Value = 1917.5
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 1917.5
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -1383.5
$ws.Range("N36").Value = -11068
$ws.Range("H54").Value = 6466.3335
$ws.Range("I54").Value = 699.5
$ws.Range("J54").Value = 18000
$ws.Range("K54").Value = 699.5
$ws.Range("L54").Value = 18000
$ws.Range("M54").Value = -215.5
$ws.Range("N54").Value = -18968
$ws.Range("H64").Value = 1592
$ws.Range("I64").Value = 1996.5
$ws.Range("J64").Value = 1389.75
$ws.Range("K64").Value = 1996.5
$ws.Range("L64").Value = 1389.75
$ws.Range("M64").Value = -1771.5
$ws.Range("N64").Value = -1839.75
$ws.Range("H67").Value = 1592
$ws.Range("I67").Value = 1996.5
$ws.Range("J67").Value = 1389.75
$ws.Range("K67").Value = 1996.5
$ws.Range("L67").Value = 1389.75
$ws.Range("M67").Value = -1216.5
$ws.Range("N67").Value = -2949.75
$ws.Range("H75").Value = 2962.6667
$ws.Range("I75").Value = 2962.6667
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2962.6667
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -2026.6667
$ws.Range("H78").Value = 2962.6667
$ws.Range("I78").Value = 2962.6667
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 8888.000100000001
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -4208.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2384.2222
$ws.Range("I16").Value = 1512.6
$ws.Range("J16").Value = 3473.75
$ws.Range("K16").Value = 1512.6
$ws.Range("L16").Value = 3473.75
$ws.Range("M16").Value = -1225.6
$ws.Range("H31").Value = 40185.75
$ws.Range("I31").Value = 1514.7273
$ws.Range("J31").Value = 72907.38
$ws.Range("K31").Value = 1514.7273
$ws.Range("L31").Value = 72907.38
$ws.Range("M31").Value = -1219.7273
$ws.Range("H32").Value = 7481
$ws.Range("I32").Value = 2605
$ws.Range("J32").Value = 9919
$ws.Range("K32").Value = 2605
$ws.Range("L32").Value = 9919
$ws.Range("M32").Value = -2289
$ws.Range("H34").Value = 40185.75
$ws.Range("I34").Value = 1514.7273
$ws.Range("J34").Value = 72907.38
$ws.Range("K34").Value = 1514.7273
$ws.Range("L34").Value = 72907.38
$ws.Range("M34").Value = -1312.7273
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H88").Value = 29797
$ws.Range("I88").Value = 10999
$ws.Range("J88").Value = 34496.5
$ws.Range("K88").Value = 10999
$ws.Range("L88").Value = 34496.5
$ws.Range("M88").Value = -10593
$ws.Range("N88").Value = -35308.5
$ws.Range("H91").Value = 29797
$ws.Range("I91").Value = 10999
$ws.Range("J91").Value = 34496.5
$ws.Range("K91").Value = 10999
$ws.Range("L91").Value = 34496.5
$ws.Range("M91").Value = -9595
$ws.Range("N91").Value = -37304.5
$ws.Range("H97").Value = 25000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 25000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 25000
$ws.Range("N97").Value = -26982
$ws.Range("H113").Value = 2384.2222
$ws.Range("I113").Value = 1512.6
$ws.Range("J113").Value = 3473.75
$ws.Range("K113").Value = 1512.6
$ws.Range("L113").Value = 3473.75
$ws.Range("M113").Value = 657.4000000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 40507.09
$ws.Range("I12").Value = 148177.5
$ws.Range("J12").Value = 130.6875
$ws.Range("K12").Value = 444532.5
$ws.Range("L12").Value = 392.0625
$ws.Range("M12").Value = -444359.5
$ws.Range("N12").Value = -738.0625
$ws.Range("H132").Value = 2734
$ws.Range("I132").Value = 2079.8
$ws.Range("J132").Value = 3097.4443
$ws.Range("K132").Value = 18718.2
$ws.Range("L132").Value = 27876.9987
$ws.Range("M132").Value = -16188.2
$ws.Range("N132").Value = -32936.9987
$ws.Range("H134").Value = 1399.76
$ws.Range("I134").Value = 1360.2174
$ws.Range("J134").Value = 1854.5
$ws.Range("K134").Value = 4080.6522
$ws.Range("L134").Value = 5563.5
$ws.Range("M134").Value = 989.3478
$ws.Range("N134").Value = -15703.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 658.3333
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 658.3333
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 658.3333
$ws.Range("N4").Value = -882.3333
$ws.Range("H32").Value = 48666
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 48666
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 48666
$ws.Range("N32").Value = -49258
$ws.Range("H46").Value = 21374.625
$ws.Range("I46").Value = 13333.333
$ws.Range("J46").Value = 26199.4
$ws.Range("K46").Value = 13333.333
$ws.Range("L46").Value = 26199.4
$ws.Range("M46").Value = -13177.333
$ws.Range("H102").Value = 7156666
$ws.Range("I102").Value = 18519876
$ws.Range("J102").Value = 1912107.9
$ws.Range("K102").Value = 18519876
$ws.Range("L102").Value = 1912107.9
$ws.Range("M102").Value = -18518254
$ws.Range("N102").Value = -1915351.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 179013.6
$ws.Range("I22").Value = 297222.66
$ws.Range("J22").Value = 1700
$ws.Range("K22").Value = 297222.66
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = -296927.66
$ws.Range("N22").Value = -2290
$ws.Range("H27").Value = 179013.6
$ws.Range("I27").Value = 297222.66
$ws.Range("J27").Value = 1700
$ws.Range("K27").Value = 297222.66
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = -297115.66
$ws.Range("N27").Value = -1914
$ws.Range("H41").Value = 28798.8
$ws.Range("I41").Value = 17998.334
$ws.Range("J41").Value = 44999.5
$ws.Range("K41").Value = 17998.334
$ws.Range("L41").Value = 44999.5
$ws.Range("M41").Value = -17560.334
$ws.Range("H46").Value = 4517.7646
$ws.Range("I46").Value = 3150
$ws.Range("J46").Value = 4938.615
$ws.Range("K46").Value = 3150
$ws.Range("L46").Value = 4938.615
$ws.Range("M46").Value = -2962
$ws.Range("N46").Value = -5314.615
$ws.Range("H53").Value = 18333
$ws.Range("I53").Value = 17499.5
$ws.Range("J53").Value = 20000
$ws.Range("K53").Value = 17499.5
$ws.Range("L53").Value = 20000
$ws.Range("M53").Value = -16981.5
$ws.Range("H68").Value = 2628.2
$ws.Range("I68").Value = 1663.5714
$ws.Range("J68").Value = 3472.25
$ws.Range("K68").Value = 1663.5714
$ws.Range("L68").Value = 3472.25
$ws.Range("M68").Value = -914.5714
$ws.Range("H71").Value = 2628.2
$ws.Range("I71").Value = 1663.5714
$ws.Range("J71").Value = 3472.25
$ws.Range("K71").Value = 8317.857
$ws.Range("L71").Value = 17361.25
$ws.Range("M71").Value = -4573.857
$ws.Range("H122").Value = 4448.3667
$ws.Range("I122").Value = 2762.8823
$ws.Range("J122").Value = 6652.4614
$ws.Range("K122").Value = 8288.6469
$ws.Range("L122").Value = 19957.3842
$ws.Range("M122").Value = -5838.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 175749.33
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 175749.33
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 175749.33
$ws.Range("N41").Value = -176529.33
$ws.Range("H113").Value = 530.55554
$ws.Range("I113").Value = 468.46155
$ws.Range("J113").Value = 615.5263
$ws.Range("K113").Value = 1405.38465
$ws.Range("L113").Value = 1846.5789
$ws.Range("M113").Value = 764.61535
$ws.Range("N113").Value = -6186.5789
$ws.Range("H133").Value = 74997.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 74997.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 74997.5
$ws.Range("N133").Value = -85117.5
